# Apply "Add data for 2022-06-01" update:
#  - Rename the sheet / update its "through" date from 2022-05-23 to 2022-05-24
#  - Update the header cell I1 text to match
#  - Update May 2022 value (I6) from 83 to 89
#  - Update Total 2022 value (I14) from 635 to 641

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Through 2022-05-24"

$ws.Range("I1").Value = "2022 (through 05-24)"
$ws.Range("I6").Value = 89
$ws.Range("I14").Value = 641
